# Daily attendance processing - 2026-01-10 20:35:58
# Swap the order of "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# in the "Recorded By" column (G) for every row where it currently appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
